$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Xpath strings in column B to remove the "[1]" index from Table[1]
$ws.Range("B2").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table/CITY"
$ws.Range("B3").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table/STATE"
$ws.Range("B4").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table/ZIP"
$ws.Range("B5").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table/AREA_CODE"
$ws.Range("B6").Value = "/Envelope/Body/GetInfoByZIPResponse/GetInfoByZIPResult/NewDataSet/Table/TIME_ZONE"

# Update the active cell selection from E8 to E6
$ws.Range("E6").Select()

$wb.Save()
